$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: D-column values that look like plain decimal numbers (e.g. "181.60")
# are written with a leading quote-prefix so Excel keeps them as literal text
# (matching the source inlineStr cells) instead of auto-converting them to
# numbers and silently dropping significant trailing zeros.

$ws.Range("D2").Value = "65.323.64"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").Value = "3.284.52"
$ws.Range("E3").Value = "  -0.91%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'576.89"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").Value = "'181.60"
$ws.Range("E6").Value = "  -3.65%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.279.17"
$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("E9").Value = "  -3.57%  "

$ws.Range("E10").Value = "  -6.40%  "

$ws.Range("D11").Value = "'0.565"
$ws.Range("E11").Value = "  -3.73%  "

$ws.Range("D12").Value = "'46.12"
$ws.Range("E12").Value = "  -3.13%  "

$ws.Range("E13").Value = "  -3.95%  "

$ws.Range("D14").Value = "'623.64"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "3.799.31"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("E16").Value = "  -3.91%  "

$ws.Range("D17").Value = "65.493.68"
$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D21").Value = "'10.84"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "'0.882"
$ws.Range("E22").Value = "  -2.97%  "

$ws.Range("D23").Value = "'17.99"
$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("D24").Value = "'99.24"
$ws.Range("E24").Value = "  -3.63%  "

$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").Value = "'3.93"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("D28").Value = "'9.30"
$ws.Range("E28").Value = "  -3.47%  "

$ws.Range("D29").Value = "'30.51"
$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("D30").Value = "'8.30"
$ws.Range("E30").Value = "  -4.66%  "

$ws.Range("D31").Value = "'6.43"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D34").Value = "'10.76"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("D35").Value = "3.800.22"
$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("D36").Value = "'0.103"
$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "'55.74"
$ws.Range("E38").Value = "  -3.91%  "

$ws.Range("E39").Value = "  -2.36%  "

$ws.Range("E40").Value = "  +6.05%  "

$ws.Range("D41").Value = "'32.17"
$ws.Range("E41").Value = "  -5.85%  "

$ws.Range("D42").Value = "'3.11"
$ws.Range("E42").Value = "  -6.28%  "

$ws.Range("E43").Value = "  -6.21%  "

$ws.Range("D44").Value = "0.0₃0670"
$ws.Range("E44").Value = "  -8.55%  "

$ws.Range("D45").Value = "'0.328"
$ws.Range("E45").Value = "  -2.84%  "

$ws.Range("D46").Value = "'0.0402"
$ws.Range("E46").Value = "  -4.23%  "

$ws.Range("D47").Value = "'3.04"
$ws.Range("E47").Value = "  -5.26%  "

$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("E49").Value = "  -3.20%  "

$ws.Range("E50").Value = "  -4.37%  "

$ws.Range("D51").Value = "'127.85"
$ws.Range("E51").Value = "  +4.40%  "

# Rows 19 & 20 swap coin identity (Chainlink <-> WrappedEther)
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.282.42"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'17.57"
$ws.Range("E20").Value = "  -2.82%  "

# Rows 32 & 33 swap coin identity (Bittensor <-> dogwifhat)
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "'3.62"
$ws.Range("E32").Value = "  -9.87%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'549.10"
$ws.Range("E33").Value = "  -1.81%  "
